$d = $word.ActiveDocument

$replacements = @(
    @("Täydellisesti", "Design: Täydellisesti"),
    @("Aikataulu", "Design: Aikataulu"),
    @("Hyvät alihankkijat", "Design: Hyvät alihankkijat"),
    @("Ei missään vaiheessa", "Design: Ei missään vaiheessa"),
    @("Pelkkää voittoa", "Design: Pelkkää voittoa"),
    @("Vähemmän virheitä kuvissa", "Design: Vähemmän virheitä kuvissa"),
    @(":)", "Design: :)")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
